$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove retired test cases that no longer apply (Repo-7, Repo-8, Result-1, Service-1)
# Delete bottom-up so earlier row numbers stay valid as we go
$ws.Cells.Item(30, 1).EntireRow.Delete()  # Service-1
$ws.Cells.Item(27, 1).EntireRow.Delete()  # Result-1
$ws.Cells.Item(26, 1).EntireRow.Delete()  # Repo-8
$ws.Cells.Item(25, 1).EntireRow.Delete()  # Repo-7

# Fill in the Description/Criteria columns for test cases that previously only had an Id
$ws.Range("B15").Value = "Getting all messages"
$ws.Range("C15").Value = "It should return 204 if nothing is found else 200"

$ws.Range("B16").Value = "Getting a single message"
$ws.Range("C16").Value = "it should return 200 if found else 404 with error"

$ws.Range("B17").Value = "Getting all lifeforms"
$ws.Range("C17").Value = "It should return 204 if nothing is found else 200"

$ws.Range("B18").Value = "posting a message"
$ws.Range("C18").Value = "it should return 200 if successful else 400 with errors"

$ws.Range("B19").Value = "Getting a single lifeform"
$ws.Range("C19").Value = "it should return 200 if found else 404 with error"

$ws.Range("B23").Value = "Can query using CQRS"
$ws.Range("C23").Value = "Entities are correctly mapped"

$ws.Range("B24").Value = "Unit of Work successfully saves all context"
$ws.Range("C24").Value = "All contexts are saved at the same time"

$ws.Range("B25").Value = "Correctly result to status conversion"
$ws.Range("C25").Value = "The Result enum is mapped to the correct http status code"

$ws.Range("B26").Value = "Hashing and salting password"
$ws.Range("C26").Value = "Same passwords always generate a different result"

$ws.Range("B27").Value = "Service method for each endpoint"
$ws.Range("C27").Value = "all methods are working correctly "

$ws.Range("B28").Value = "Service method for each endpoint"
$ws.Range("C28").Value = "all methods are working correctly "

$ws.Range("B29").Value = "Service method for each endpoint"
$ws.Range("C29").Value = "all methods are working correctly "

$ws.Range("B30").Value = "All validations for creation are working"
$ws.Range("C30").Value = "If err they return error messages"

$ws.Range("B31").Value = "All validations for queying work"
$ws.Range("C31").Value = "the querying understand them and the correct data is found"

# Re-apply the existing sort (by Id) so the sheet's sort-state range reflects the new row count
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A31"))
$ws.Sort.SetRange($ws.Range("A2:E31"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Restore the author's last-known selection
$ws.Range("D13").Select()
